$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Apio (Macroferia Regional de Talca).
# It belongs right after the existing row 222, so insert a fresh row at 223
# which pushes the former rows 223-241 down to 224-242 (carrying their
# original data/format with them - this mirrors the diff, where every row
# from 224 to 242 now holds what used to be one row above it).
$ws.Rows(223).Insert()

# Fill in the new record's data in row 223.
$ws.Range("A223").Value = 5
$ws.Range("B223").Value = "Macroferia Regional de Talca"
$ws.Range("C223").Value = "Maule"
$ws.Range("D223").Value = 44826
$ws.Range("E223").Value = 7
$ws.Range("F223").Value = 100112017
$ws.Range("G223").Value = "Apio"
$ws.Range("H223").Value = "Americana (o)"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 700
$ws.Range("K223").Value = 9000
$ws.Range("L223").Value = 9000
$ws.Range("M223").Value = 9000
$ws.Range("N223").Value = "$/docena de matas"
$ws.Range("O223").Value = "Provincia del Elquí"
$ws.Range("P223").Value = 1500
$ws.Range("Q223").Value = 6
$ws.Range("R223").Value = "Hortaliza"
